$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "duplicate_image_filename" column (E) needs "NA" filled in for every
# data row of the main stimuli table (rows 2-21), matching the header in E1.
$ws.Range("E2:E21").Value = "NA"
